# Generate Report for Handoff
# Adds a new row (row 3) to each of the Overview / zh-cn / de-de sheets
# representing the new file dad3e217-...md that is now "Ready for handoff".

$wb = $excel.ActiveWorkbook

$mdA = "dad3e217-0687-4807-9617-f4259766d8ecoooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$mdBdisp = "e2e\dad3e217-0687-4807-9617-f4259766d8ecoooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$status = "Ready for handoff"
$date154946 = "2016-09-07 15:39:46"
$xlfZhCn = "dad3e217-0687-4807-9617-f4259766d8ecooooooooooooooooooooooooooooooooooooooo.b4afa2c6d7328817886ee3d0e657d5c5d16af3df.zh-cn.xlf"
$date153935 = "2016-09-07 15:39:35"
$xlfDeDe = "dad3e217-0687-4807-9617-f4259766d8ecooooooooooooooooooooooooooooooooooooooo.b4afa2c6d7328817886ee3d0e657d5c5d16af3df.de-de.xlf"
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c378a3f0142533db4d0e53dada88f7104d1d6214/e2e/dad3e217-0687-4807-9617-f4259766d8ecoooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$date0001 = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdA
$wsOverview.Range("B3").Value = $mdBdisp
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, "", "", $mdBdisp) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $date154946
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $mdA
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkTarget, "", "", $mdA) | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $date153935
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = $date0001
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $mdA
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkTarget, "", "", $mdA) | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $date154946
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = $date0001
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3

Write-Host "Report generated for handoff."
